$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("en")

# The "Metabolism" category becomes "Cell Structure" (in place), the old
# "Coating" category row is replaced by "Motility" (pulled up from row 8),
# and the rows below shift up accordingly.
$ws.Cells.Item(6, 1).Value = "categoryCellStructure"
$ws.Cells.Item(6, 2).Value = "Cell Structure"

$ws.Cells.Item(7, 1).Value = "categoryMotility"
$ws.Cells.Item(7, 2).Value = "Motility"

$ws.Cells.Item(8, 1).Value = "nucleoid"
$ws.Cells.Item(8, 2).Value = "Nucleoid"

$ws.Cells.Item(9, 1).Value = "ribosome"
$ws.Cells.Item(9, 2).Value = "Ribosome"

$ws.Cells.Item(10, 1).Value = "test1"
$ws.Cells.Item(10, 2).Value = "Test 1"

$ws.Cells.Item(11, 1).Value = "test2"
$ws.Cells.Item(11, 2).Value = "Test 2"

$ws.Cells.Item(12, 1).Value = "testBodyCapsule"
$ws.Cells.Item(12, 2).Value = "Capsule"

$ws.Cells.Item(13, 1).Value = "testBodySphere"
$ws.Cells.Item(13, 2).Value = "Sphere"

# Four new extremophile attribute rows for the energy/hazard display info.
# Row 14's value was typed before its key, so the underlying shared-string
# table picks up "Thermophile" ahead of "cellStructureThermophile".
$ws.Cells.Item(14, 2).Value = "Thermophile"
$ws.Cells.Item(14, 1).Value = "cellStructureThermophile"

$ws.Cells.Item(15, 1).Value = "cellStructurePsychrophile"
$ws.Cells.Item(15, 2).Value = "Psychrophile"

$ws.Cells.Item(16, 1).Value = "cellStructureMethanogen"
$ws.Cells.Item(16, 2).Value = "Methanogen"

$ws.Cells.Item(17, 1).Value = "cellStructureHalophile"
$ws.Cells.Item(17, 2).Value = "Halophile"

$ws.Range("B17").Select()
